$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.836.98'
$ws.Range("E2").Value = '  -0.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.584.67'
$ws.Range("E3").Value = '  +2.34%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.93'
$ws.Range("E5").Value = '  +1.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.02'
$ws.Range("E6").Value = '  +4.01%  '

$ws.Range("E7").Value = '  +0.86%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.550'
$ws.Range("E9").Value = '  +0.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.80'
$ws.Range("E10").Value = '  +1.56%  '

$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.72'
$ws.Range("E12").Value = '  +2.06%  '

$ws.Range("E13").Value = '  +6.96%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.573.04'
$ws.Range("E14").Value = '  +1.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.886'
$ws.Range("E15").Value = '  +2.40%  '

$ws.Range("E16").Value = '  +2.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '42.918.60'
$ws.Range("E17").Value = '  +0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.88'
$ws.Range("E18").Value = '  +4.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0990'
$ws.Range("E19").Value = '  +2.89%  '

$ws.Range("E20").Value = '  +1.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.95'
$ws.Range("E21").Value = '  -1.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '254.57'
$ws.Range("E22").Value = '  -1.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.96'
$ws.Range("E23").Value = '  +2.24%  '

$ws.Range("E24").Value = '  -1.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '28.63'
$ws.Range("E25").Value = '  -0.97%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.21'
$ws.Range("E27").Value = '  +2.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.98'
$ws.Range("E28").Value = '  +6.01%  '

$ws.Range("E29").Value = '  -0.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.04'
$ws.Range("E30").Value = '  +1.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.58'
$ws.Range("E31").Value = '  +2.77%  '

$ws.Range("E32").Value = '  -0.84%  '

$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0815'
$ws.Range("E34").Value = '  +1.98%  '

$ws.Range("E35").Value = '  -3.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.33'
$ws.Range("E36").Value = '  +12.81%  '

$ws.Range("E37").Value = '  +0.41%  '

$ws.Range("E38").Value = '  +1.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.38'
$ws.Range("E39").Value = '  -1.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.09'
$ws.Range("E40").Value = '  +30.43%  '

$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.42'
$ws.Range("E41").Value = '  -0.38%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0311'
$ws.Range("E42").Value = '  +0.95%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.88'
$ws.Range("E43").Value = '  +1.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.071.57'
$ws.Range("E44").Value = '  +2.74%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.26'
$ws.Range("E46").Value = '  +4.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.23'
$ws.Range("E47").Value = '  -0.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.32'
$ws.Range("E48").Value = '  +11.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.836.19'
$ws.Range("E49").Value = '  +2.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '106.18'
$ws.Range("E50").Value = '  +3.29%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.67'
$ws.Range("E51").Value = '  +2.62%  '
